# Commit: update 2548 files and delete 4 files
# This sheet: re-order monthly rows so each calendar year lists Oct, Nov, Dec
# first followed by Jan-Sep (matches the source publication order), and append
# newly released data through 2022 and partial 2023 (columns B/C only; D has not
# been published yet for these months).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 50:68 are brand new -- prime column A with the same style (s="1": bold,
# boxed border, centered) used by the existing date column before writing values.
$ws.Range("A49").Copy()
$ws.Range("A50:A68").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(2, 1).Value = "2018-10"
$ws.Cells.Item(2, 2).Value = 101.9
$ws.Cells.Item(2, 3).Value = 100.5
$ws.Cells.Item(2, 4).Value = 99.2

$ws.Cells.Item(3, 1).Value = "2018-11"
$ws.Cells.Item(3, 2).Value = 101.7
$ws.Cells.Item(3, 3).Value = 100.6
$ws.Cells.Item(3, 4).Value = 99.7

$ws.Cells.Item(4, 1).Value = "2018-12"
$ws.Cells.Item(4, 2).Value = 101.5
$ws.Cells.Item(4, 3).Value = 100.8
$ws.Cells.Item(4, 4).Value = 99.2

$ws.Cells.Item(5, 1).Value = "2018-01"
$ws.Cells.Item(5, 2).Value = 101.7488
$ws.Cells.Item(5, 3).Value = 102.8955
$ws.Cells.Item(5, 4).Value = 98.7358

$ws.Cells.Item(6, 1).Value = "2018-02"
$ws.Cells.Item(6, 2).Value = 101.8
$ws.Cells.Item(6, 3).Value = 102.6
$ws.Cells.Item(6, 4).Value = 99.1

$ws.Cells.Item(7, 1).Value = "2018-03"
$ws.Cells.Item(7, 2).Value = 101.7
$ws.Cells.Item(7, 3).Value = 102.8
$ws.Cells.Item(7, 4).Value = 97.6

$ws.Cells.Item(8, 1).Value = "2018-04"
$ws.Cells.Item(8, 2).Value = 102
$ws.Cells.Item(8, 3).Value = 102.8
$ws.Cells.Item(8, 4).Value = 98.1

$ws.Cells.Item(9, 1).Value = "2018-05"
$ws.Cells.Item(9, 2).Value = 102
$ws.Cells.Item(9, 3).Value = 102.4
$ws.Cells.Item(9, 4).Value = 98.4

$ws.Cells.Item(10, 1).Value = "2018-06"
$ws.Cells.Item(10, 2).Value = 101.9
$ws.Cells.Item(10, 3).Value = 101.2
$ws.Cells.Item(10, 4).Value = 98.4

$ws.Cells.Item(11, 1).Value = "2018-07"
$ws.Cells.Item(11, 2).Value = 101.7
$ws.Cells.Item(11, 3).Value = 101.6
$ws.Cells.Item(11, 4).Value = 98.7

$ws.Cells.Item(12, 1).Value = "2018-08"
$ws.Cells.Item(12, 2).Value = 101.7
$ws.Cells.Item(12, 3).Value = 101.5
$ws.Cells.Item(12, 4).Value = 98.5

$ws.Cells.Item(13, 1).Value = "2018-09"
$ws.Cells.Item(13, 2).Value = 101.6
$ws.Cells.Item(13, 3).Value = 101.1
$ws.Cells.Item(13, 4).Value = 98.5

$ws.Cells.Item(14, 1).Value = "2019-10"
$ws.Cells.Item(14, 2).Value = 99.2
$ws.Cells.Item(14, 3).Value = 101.3
$ws.Cells.Item(14, 4).Value = 98.7

$ws.Cells.Item(15, 1).Value = "2019-11"
$ws.Cells.Item(15, 2).Value = 99.3
$ws.Cells.Item(15, 3).Value = 100.7
$ws.Cells.Item(15, 4).Value = 98.8

$ws.Cells.Item(16, 1).Value = "2019-12"
$ws.Cells.Item(16, 2).Value = 99.3
$ws.Cells.Item(16, 3).Value = 100.8
$ws.Cells.Item(16, 4).Value = 99.1

$ws.Cells.Item(17, 1).Value = "2019-01"
$ws.Cells.Item(17, 2).Value = 101.2
$ws.Cells.Item(17, 3).Value = 101
$ws.Cells.Item(17, 4).Value = 98.9

$ws.Cells.Item(18, 1).Value = "2019-02"
$ws.Cells.Item(18, 2).Value = 100.9
$ws.Cells.Item(18, 3).Value = 101.3
$ws.Cells.Item(18, 4).Value = 99.6

$ws.Cells.Item(19, 1).Value = "2019-03"
$ws.Cells.Item(19, 2).Value = 100.7
$ws.Cells.Item(19, 3).Value = 101.1
$ws.Cells.Item(19, 4).Value = 99.1

$ws.Cells.Item(20, 1).Value = "2019-04"
$ws.Cells.Item(20, 2).Value = 100.6
$ws.Cells.Item(20, 3).Value = 101.4
$ws.Cells.Item(20, 4).Value = 99.7

$ws.Cells.Item(21, 1).Value = "2019-05"
$ws.Cells.Item(21, 2).Value = 100.4
$ws.Cells.Item(21, 3).Value = 102.2
$ws.Cells.Item(21, 4).Value = 100.5

$ws.Cells.Item(22, 1).Value = "2019-06"
$ws.Cells.Item(22, 2).Value = 100.3
$ws.Cells.Item(22, 3).Value = 102.6
$ws.Cells.Item(22, 4).Value = 100.2

$ws.Cells.Item(23, 1).Value = "2019-07"
$ws.Cells.Item(23, 2).Value = 100.5
$ws.Cells.Item(23, 3).Value = 101.8
$ws.Cells.Item(23, 4).Value = 98

$ws.Cells.Item(24, 1).Value = "2019-08"
$ws.Cells.Item(24, 2).Value = 100.3
$ws.Cells.Item(24, 3).Value = 102.1
$ws.Cells.Item(24, 4).Value = 99.2

$ws.Cells.Item(25, 1).Value = "2019-09"
$ws.Cells.Item(25, 2).Value = 100
$ws.Cells.Item(25, 3).Value = 101.9
$ws.Cells.Item(25, 4).Value = 99

$ws.Cells.Item(26, 1).Value = "2020-10"
$ws.Cells.Item(26, 2).Value = 98.2
$ws.Cells.Item(26, 3).Value = 98.9
$ws.Cells.Item(26, 4).Value = 98.7

$ws.Cells.Item(27, 1).Value = "2020-11"
$ws.Cells.Item(27, 2).Value = 97.9
$ws.Cells.Item(27, 3).Value = 98.7
$ws.Cells.Item(27, 4).Value = 98.3

$ws.Cells.Item(28, 1).Value = "2020-12"
$ws.Cells.Item(28, 2).Value = 97.9
$ws.Cells.Item(28, 3).Value = 98.5
$ws.Cells.Item(28, 4).Value = 98.5

$ws.Cells.Item(29, 1).Value = "2020-01"
$ws.Cells.Item(29, 2).Value = 99.4
$ws.Cells.Item(29, 3).Value = 100.9
$ws.Cells.Item(29, 4).Value = 99.7

$ws.Cells.Item(30, 1).Value = "2020-02"
$ws.Cells.Item(30, 2).Value = 99.4
$ws.Cells.Item(30, 3).Value = 100.8
$ws.Cells.Item(30, 4).Value = 99.5

$ws.Cells.Item(31, 1).Value = "2020-03"
$ws.Cells.Item(31, 2).Value = 99.6
$ws.Cells.Item(31, 3).Value = 100.1
$ws.Cells.Item(31, 4).Value = 99.7

$ws.Cells.Item(32, 1).Value = "2020-04"
$ws.Cells.Item(32, 2).Value = 99.3
$ws.Cells.Item(32, 3).Value = 99.5
$ws.Cells.Item(32, 4).Value = 99.4

$ws.Cells.Item(33, 1).Value = "2020-05"
$ws.Cells.Item(33, 2).Value = 99.3
$ws.Cells.Item(33, 3).Value = 99.2
$ws.Cells.Item(33, 4).Value = 98.3

$ws.Cells.Item(34, 1).Value = "2020-06"
$ws.Cells.Item(34, 2).Value = 98.6
$ws.Cells.Item(34, 3).Value = 99.3
$ws.Cells.Item(34, 4).Value = 97.2

$ws.Cells.Item(35, 1).Value = "2020-07"
$ws.Cells.Item(35, 2).Value = 98.2
$ws.Cells.Item(35, 3).Value = 100
$ws.Cells.Item(35, 4).Value = 99.2

$ws.Cells.Item(36, 1).Value = "2020-08"
$ws.Cells.Item(36, 2).Value = 98.3
$ws.Cells.Item(36, 3).Value = 99.4
$ws.Cells.Item(36, 4).Value = 98.2

$ws.Cells.Item(37, 1).Value = "2020-09"
$ws.Cells.Item(37, 2).Value = 98.3
$ws.Cells.Item(37, 3).Value = 99.4
$ws.Cells.Item(37, 4).Value = 98.9

$ws.Cells.Item(38, 1).Value = "2021-10"
$ws.Cells.Item(38, 2).Value = 101
$ws.Cells.Item(38, 3).Value = 103.2
$ws.Cells.Item(38, 4).Value = 100

$ws.Cells.Item(39, 1).Value = "2021-11"
$ws.Cells.Item(39, 2).Value = 101.5
$ws.Cells.Item(39, 3).Value = 103.7
$ws.Cells.Item(39, 4).Value = 100

$ws.Cells.Item(40, 1).Value = "2021-12"
$ws.Cells.Item(40, 2).Value = 101.5
$ws.Cells.Item(40, 3).Value = 101.2
$ws.Cells.Item(40, 4).Value = 100

$ws.Cells.Item(41, 1).Value = "2021-01"
$ws.Cells.Item(41, 2).Value = 98.9
$ws.Cells.Item(41, 3).Value = 98.4
$ws.Cells.Item(41, 4).Value = 96.8

$ws.Cells.Item(42, 1).Value = "2021-02"
$ws.Cells.Item(42, 2).Value = 99
$ws.Cells.Item(42, 3).Value = 98.9
$ws.Cells.Item(42, 4).Value = 96.8

$ws.Cells.Item(43, 1).Value = "2021-03"
$ws.Cells.Item(43, 2).Value = 99.4
$ws.Cells.Item(43, 3).Value = 100.6
$ws.Cells.Item(43, 4).Value = 96.8

$ws.Cells.Item(44, 1).Value = "2021-04"
$ws.Cells.Item(44, 2).Value = 99.7
$ws.Cells.Item(44, 3).Value = 99.5
$ws.Cells.Item(44, 4).Value = 96.8

$ws.Cells.Item(45, 1).Value = "2021-05"
$ws.Cells.Item(45, 2).Value = 100.8
$ws.Cells.Item(45, 3).Value = 99.4
$ws.Cells.Item(45, 4).Value = 100

$ws.Cells.Item(46, 1).Value = "2021-06"
$ws.Cells.Item(46, 2).Value = 100.9
$ws.Cells.Item(46, 3).Value = 99.6
$ws.Cells.Item(46, 4).Value = 100

$ws.Cells.Item(47, 1).Value = "2021-07"
$ws.Cells.Item(47, 2).Value = 101.2
$ws.Cells.Item(47, 3).Value = 100.2
$ws.Cells.Item(47, 4).Value = 100

$ws.Cells.Item(48, 1).Value = "2021-08"
$ws.Cells.Item(48, 2).Value = 100.9
$ws.Cells.Item(48, 3).Value = 102.1
$ws.Cells.Item(48, 4).Value = 100

$ws.Cells.Item(49, 1).Value = "2021-09"
$ws.Cells.Item(49, 2).Value = 100.8
$ws.Cells.Item(49, 3).Value = 103.4
$ws.Cells.Item(49, 4).Value = 100

$ws.Cells.Item(50, 1).Value = "2022-10"
$ws.Cells.Item(50, 2).Value = 100.9
$ws.Cells.Item(50, 3).Value = 96.4
$ws.Cells.Item(50, 4).ClearContents()

$ws.Cells.Item(51, 1).Value = "2022-11"
$ws.Cells.Item(51, 2).Value = 100.4
$ws.Cells.Item(51, 3).Value = 96.3
$ws.Cells.Item(51, 4).ClearContents()

$ws.Cells.Item(52, 1).Value = "2022-12"
$ws.Cells.Item(52, 2).Value = 100.1
$ws.Cells.Item(52, 3).Value = 102.6
$ws.Cells.Item(52, 4).ClearContents()

$ws.Cells.Item(53, 1).Value = "2022-01"
$ws.Cells.Item(53, 2).Value = 102
$ws.Cells.Item(53, 3).Value = 101
$ws.Cells.Item(53, 4).ClearContents()

$ws.Cells.Item(54, 1).Value = "2022-02"
$ws.Cells.Item(54, 2).Value = 101.7
$ws.Cells.Item(54, 3).Value = 100.7
$ws.Cells.Item(54, 4).ClearContents()

$ws.Cells.Item(55, 1).Value = "2022-03"
$ws.Cells.Item(55, 2).Value = 101.5
$ws.Cells.Item(55, 3).Value = 98.2
$ws.Cells.Item(55, 4).ClearContents()

$ws.Cells.Item(56, 1).Value = "2022-04"
$ws.Cells.Item(56, 2).Value = 101.4
$ws.Cells.Item(56, 3).Value = 99.3
$ws.Cells.Item(56, 4).ClearContents()

$ws.Cells.Item(57, 1).Value = "2022-05"
$ws.Cells.Item(57, 2).Value = 100.7
$ws.Cells.Item(57, 3).Value = 99.1
$ws.Cells.Item(57, 4).ClearContents()

$ws.Cells.Item(58, 1).Value = "2022-06"
$ws.Cells.Item(58, 2).Value = 100.9
$ws.Cells.Item(58, 3).Value = 98.5
$ws.Cells.Item(58, 4).ClearContents()

$ws.Cells.Item(59, 1).Value = "2022-07"
$ws.Cells.Item(59, 2).Value = 100.9
$ws.Cells.Item(59, 3).Value = 97.6
$ws.Cells.Item(59, 4).ClearContents()

$ws.Cells.Item(60, 1).Value = "2022-08"
$ws.Cells.Item(60, 2).Value = 100.8
$ws.Cells.Item(60, 3).Value = 97.2
$ws.Cells.Item(60, 4).ClearContents()

$ws.Cells.Item(61, 1).Value = "2022-09"
$ws.Cells.Item(61, 2).Value = 100.9
$ws.Cells.Item(61, 3).Value = 95.5
$ws.Cells.Item(61, 4).ClearContents()

$ws.Cells.Item(62, 1).Value = "2023-01"
$ws.Cells.Item(62, 2).Value = 99.7
$ws.Cells.Item(62, 3).Value = 101.7
$ws.Cells.Item(62, 4).ClearContents()

$ws.Cells.Item(63, 1).Value = "2023-02"
$ws.Cells.Item(63, 2).Value = 99.4
$ws.Cells.Item(63, 3).Value = 101.2
$ws.Cells.Item(63, 4).ClearContents()

$ws.Cells.Item(64, 1).Value = "2023-03"
$ws.Cells.Item(64, 2).Value = 99.4
$ws.Cells.Item(64, 3).Value = 102.8
$ws.Cells.Item(64, 4).ClearContents()

$ws.Cells.Item(65, 1).Value = "2023-04"
$ws.Cells.Item(65, 2).Value = 99.2
$ws.Cells.Item(65, 3).Value = 102.8
$ws.Cells.Item(65, 4).ClearContents()

$ws.Cells.Item(66, 1).Value = "2023-05"
$ws.Cells.Item(66, 2).Value = 99
$ws.Cells.Item(66, 3).Value = 145.3
$ws.Cells.Item(66, 4).ClearContents()

$ws.Cells.Item(67, 1).Value = "2023-06"
$ws.Cells.Item(67, 2).Value = 98.9
$ws.Cells.Item(67, 3).Value = 145.8
$ws.Cells.Item(67, 4).ClearContents()

$ws.Cells.Item(68, 1).Value = "2023-07"
$ws.Cells.Item(68, 2).Value = 98.4
$ws.Cells.Item(68, 3).Value = 146.5
$ws.Cells.Item(68, 4).ClearContents()

"done"
